$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the differing species-record data between rows 6-10
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Snapshot the "before" values of rows 6-10 for the varying columns
$data = @{}
foreach ($r in 6..10) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value()
    }
    $data[$r] = $rowVals
}

# Row 7 currently has an (empty) L cell that no other row in 6-10 has.
# After the edit that empty L cell needs to move along with row 7's record
# to row 9 (since row 9 will now hold what used to be row 7's data).
# Do this via a genuine cell copy so the empty cell is preserved faithfully,
# then blank out the source so it disappears from row 7.
$ws.Range("L7").Copy($ws.Range("L9"))
$ws.Range("L7").Value = ""

# Mapping: new row -> old row whose record it should now hold
$mapping = @{ 6 = 10; 7 = 8; 8 = 9; 9 = 7; 10 = 6 }

foreach ($newRow in 6..10) {
    $oldRow = $mapping[$newRow]
    $src = $data[$oldRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value = $src[$col]
    }
}
